# Reorder each year's 12 monthly rows so Oct/Nov/Dec move to the
# front of their year block (ahead of Jan-Sep), matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "2014-10"
$ws.Cells.Item(2, 2).Value = 86.2209
$ws.Cells.Item(2, 3).Value = 93.0328
$ws.Cells.Item(3, 1).Value = "2014-11"
$ws.Cells.Item(3, 2).Value = 85.3546
$ws.Cells.Item(3, 3).Value = 89.4446
$ws.Cells.Item(4, 1).Value = "2014-12"
$ws.Cells.Item(4, 2).Value = 84.6069
$ws.Cells.Item(4, 3).Value = 83.5389
$ws.Cells.Item(5, 1).Value = "2014-01"
$ws.Cells.Item(5, 2).Value = 91.1606
$ws.Cells.Item(5, 3).Value = 98.9873
$ws.Cells.Item(6, 1).Value = "2014-02"
$ws.Cells.Item(6, 2).Value = 88.3893
$ws.Cells.Item(6, 3).Value = 96.874
$ws.Cells.Item(7, 1).Value = "2014-03"
$ws.Cells.Item(7, 2).Value = 85.8698
$ws.Cells.Item(7, 3).Value = 95.9156
$ws.Cells.Item(8, 1).Value = "2014-04"
$ws.Cells.Item(8, 2).Value = 85.0281
$ws.Cells.Item(8, 3).Value = 98.0948
$ws.Cells.Item(9, 1).Value = "2014-05"
$ws.Cells.Item(9, 2).Value = 85.6277
$ws.Cells.Item(9, 3).Value = 101.787
$ws.Cells.Item(10, 1).Value = "2014-06"
$ws.Cells.Item(10, 2).Value = 86.5925
$ws.Cells.Item(10, 3).Value = 102.8929
$ws.Cells.Item(11, 1).Value = "2014-07"
$ws.Cells.Item(11, 2).Value = 88.1713
$ws.Cells.Item(11, 3).Value = 103.1263
$ws.Cells.Item(12, 1).Value = "2014-08"
$ws.Cells.Item(12, 2).Value = 87.984
$ws.Cells.Item(12, 3).Value = 99.4811
$ws.Cells.Item(13, 1).Value = "2014-09"
$ws.Cells.Item(13, 2).Value = 86.8655
$ws.Cells.Item(13, 3).Value = 94.9137
$ws.Cells.Item(14, 1).Value = "2015-10"
$ws.Cells.Item(14, 2).Value = 83.8
$ws.Cells.Item(14, 3).Value = 76.3
$ws.Cells.Item(15, 1).Value = "2015-11"
$ws.Cells.Item(15, 2).Value = 82.4802
$ws.Cells.Item(15, 3).Value = 78.8131
$ws.Cells.Item(16, 1).Value = "2015-12"
$ws.Cells.Item(16, 2).Value = 80.6932
$ws.Cells.Item(16, 3).Value = 81.3818
$ws.Cells.Item(17, 1).Value = "2015-01"
$ws.Cells.Item(17, 2).Value = 84.1273
$ws.Cells.Item(17, 3).Value = 75.7457
$ws.Cells.Item(18, 1).Value = "2015-02"
$ws.Cells.Item(18, 2).Value = 85.1134
$ws.Cells.Item(18, 3).Value = 74.1905
$ws.Cells.Item(19, 1).Value = "2015-03"
$ws.Cells.Item(19, 2).Value = 86.4218
$ws.Cells.Item(19, 3).Value = 77.337
$ws.Cells.Item(20, 1).Value = "2015-04"
$ws.Cells.Item(20, 2).Value = 87.8343
$ws.Cells.Item(20, 3).Value = 77.0067
$ws.Cells.Item(21, 1).Value = "2015-05"
$ws.Cells.Item(21, 2).Value = 87.4366
$ws.Cells.Item(21, 3).Value = 80.4051
$ws.Cells.Item(22, 1).Value = "2015-06"
$ws.Cells.Item(22, 2).Value = 87.2244
$ws.Cells.Item(22, 3).Value = 80.247
$ws.Cells.Item(23, 1).Value = "2015-07"
$ws.Cells.Item(23, 2).Value = 86.6889
$ws.Cells.Item(23, 3).Value = 77.6766
$ws.Cells.Item(24, 1).Value = "2015-08"
$ws.Cells.Item(24, 2).Value = 85.9906
$ws.Cells.Item(24, 3).Value = 73.8117
$ws.Cells.Item(25, 1).Value = "2015-09"
$ws.Cells.Item(25, 2).Value = 84.7564
$ws.Cells.Item(25, 3).Value = 73.7792
$ws.Cells.Item(26, 1).Value = "2016-10"
$ws.Cells.Item(26, 2).Value = 130.2
$ws.Cells.Item(26, 3).Value = 99.4
$ws.Cells.Item(27, 1).Value = "2016-11"
$ws.Cells.Item(27, 2).Value = 151.9
$ws.Cells.Item(27, 3).Value = 102.4
$ws.Cells.Item(28, 1).Value = "2016-12"
$ws.Cells.Item(28, 2).Value = 164.8
$ws.Cells.Item(28, 3).Value = 109.1
$ws.Cells.Item(29, 1).Value = "2016-01"
$ws.Cells.Item(29, 2).Value = 79.1724
$ws.Cells.Item(29, 3).Value = 87.7713
$ws.Cells.Item(30, 1).Value = "2016-02"
$ws.Cells.Item(30, 2).Value = 79.1816
$ws.Cells.Item(30, 3).Value = 86.8593
$ws.Cells.Item(31, 1).Value = "2016-03"
$ws.Cells.Item(31, 2).Value = 81.094
$ws.Cells.Item(31, 3).Value = 82.1733
$ws.Cells.Item(32, 1).Value = "2016-04"
$ws.Cells.Item(32, 2).Value = 84.6932
$ws.Cells.Item(32, 3).Value = 83.8197
$ws.Cells.Item(33, 1).Value = "2016-05"
$ws.Cells.Item(33, 2).Value = 93.2
$ws.Cells.Item(33, 3).Value = 83.0
$ws.Cells.Item(34, 1).Value = "2016-06"
$ws.Cells.Item(34, 2).Value = 99.4
$ws.Cells.Item(34, 3).Value = 86.5
$ws.Cells.Item(35, 1).Value = "2016-07"
$ws.Cells.Item(35, 2).Value = 100.8
$ws.Cells.Item(35, 3).Value = 89.3
$ws.Cells.Item(36, 1).Value = "2016-08"
$ws.Cells.Item(36, 2).Value = 106.4
$ws.Cells.Item(36, 3).Value = 92.9
$ws.Cells.Item(37, 1).Value = "2016-09"
$ws.Cells.Item(37, 2).Value = 116.7
$ws.Cells.Item(37, 3).Value = 97.3
$ws.Cells.Item(38, 1).Value = "2017-10"
$ws.Cells.Item(38, 2).Value = 142.1
$ws.Cells.Item(38, 3).Value = 110.5
$ws.Cells.Item(39, 1).Value = "2017-11"
$ws.Cells.Item(39, 2).Value = 116.9
$ws.Cells.Item(39, 3).Value = 113.2
$ws.Cells.Item(40, 1).Value = "2017-12"
$ws.Cells.Item(40, 2).Value = 116.6
$ws.Cells.Item(40, 3).Value = 111.5
$ws.Cells.Item(41, 1).Value = "2017-01"
$ws.Cells.Item(41, 2).Value = 170.3
$ws.Cells.Item(41, 3).Value = 117.1
$ws.Cells.Item(42, 1).Value = "2017-02"
$ws.Cells.Item(42, 2).Value = 173.4
$ws.Cells.Item(42, 3).Value = 123.4
$ws.Cells.Item(43, 1).Value = "2017-03"
$ws.Cells.Item(43, 2).Value = 170.6
$ws.Cells.Item(43, 3).Value = 123.2
$ws.Cells.Item(44, 1).Value = "2017-04"
$ws.Cells.Item(44, 2).Value = 173.0
$ws.Cells.Item(44, 3).Value = 120.1
$ws.Cells.Item(45, 1).Value = "2017-05"
$ws.Cells.Item(45, 2).Value = 160.7
$ws.Cells.Item(45, 3).Value = 115.5
$ws.Cells.Item(46, 1).Value = "2017-06"
$ws.Cells.Item(46, 2).Value = 149.6
$ws.Cells.Item(46, 3).Value = 108.9
$ws.Cells.Item(47, 1).Value = "2017-07"
$ws.Cells.Item(47, 2).Value = 152.2
$ws.Cells.Item(47, 3).Value = 104.1
$ws.Cells.Item(48, 1).Value = "2017-08"
$ws.Cells.Item(48, 2).Value = 153.7
$ws.Cells.Item(48, 3).Value = 110.3
$ws.Cells.Item(49, 1).Value = "2017-09"
$ws.Cells.Item(49, 2).Value = 153.2
$ws.Cells.Item(49, 3).Value = 109.8
